$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (rows 2-6), columns A-N
$data = @(
    @(3011, "Tata Indigo eCS eLS BS IV",         "Tata", "Kolkata",    2016, 78000,  "Diesel", "Manual", "First",  1396, 69.01000000000001, 5, 25,    2.3),
    @(4748, "Tata Tiago 1.2 Revotron XT Option", "Tata", "Coimbatore", 2019, 16698,  "Petrol", "Manual", "First",  1199, 84,                 5, 23.84, 5.85),
    @(3263, "Tata Indica Vista Quadrajet LS",     "Tata", "Kolkata",    2012, 80000,  "Diesel", "Manual", "Second", 1248, 74,                 5, 22.3,  1.2),
    @(3551, "Tata Indigo GLE",                    "Tata", "Chennai",    2008, 73000,  "Petrol", "Manual", "First",  1396, 65,                 5, 14.4,  0.8),
    @(319,  "Tata Indigo LS Dicor",                "Tata", "Pune",       2008, 138000, "Diesel", "Manual", "Second", 1396, 70,                 5, 17.1,  1.25)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
